$wb = $excel.ActiveWorkbook

# Add new shared string / comment text to MAIN sheet, new row
$mainSheet = $wb.Worksheets.Item("MAIN")
$mainSheet.Range("B7").Value = "2019.11.11 Femto에 V2.0으로 release"

# Rename MAIN -> MAIN V2.0
$mainSheet.Name = "MAIN V2.0"

# Set active sheet / selection on MAIN V2.0
$mainSheet.Activate()
$mainSheet.Range("F16").Select()

# Set selection on Issue sheet (no longer tab-selected)
$issueSheet = $wb.Worksheets.Item("Issue")
$issueSheet.Range("G24").Select()
